$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hit_miss_rule")

# Update the underlying raw values (column H, I, J) that feed the
# ROUND()/SUM() formulas in columns D, E, F. Excel will recalculate the
# dependent formula cells automatically.

$ws.Range("H5").Value = 91.224967956542969
$ws.Range("I6").Value = 8.7750320434570313

$ws.Range("H8").Value = 3.7213354110717769
$ws.Range("I8").Value = 4.5544266700744629

$ws.Range("H9").Value = 6.5950651168823242
$ws.Range("I9").Value = 7.4397678375244141

$ws.Range("H10").Value = 77.015411376953125
$ws.Range("I10").Value = 17.839620590209961
$ws.Range("J10").Value = 67.911468505859375

$excel.CalculateFullRebuild()
